$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.0292345
$ws.Range("H2").Value = 0.058469
$ws.Range("I2").Value = 0.4428765120700495
$ws.Range("J2").Value = 0.346386487911515
$ws.Range("M2").Value = 3.1851815
$ws.Range("N2").Value = 6.370363
$ws.Range("O2").Value = 0.4406530230187619
$ws.Range("P2").Value = 0.3851702893788179
$ws.Range("Q2").Value = 0.09311718856175
$ws.Range("R2").Value = 0.372468754247
$ws.Range("S2").Value = 0.1951548738676725
$ws.Range("T2").Value = 0.1334177837857906
# Row 3
$ws.Range("G3").Value = 0.0292345
$ws.Range("H3").Value = 0.058469
$ws.Range("I3").Value = 0.4428765120700495
$ws.Range("J3").Value = 0.346386487911515
$ws.Range("O3").Value = 0.2827048402157753
$ws.Range("P3").Value = 0.3706641033643825
$ws.Range("Q3").Value = 0.05974015503933333
$ws.Range("R3").Value = 0.358440930236
$ws.Range("S3").Value = 0.1252033335800832
$ws.Range("T3").Value = 0.1283930369592592
# Row 4
$ws.Range("G4").Value = 0.0292345
$ws.Range("H4").Value = 0.058469
$ws.Range("I4").Value = 0.4428765120700495
$ws.Range("J4").Value = 0.346386487911515
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.003928
$ws.Range("N4").Value = 0.011784
$ws.Range("O4").Value = 0.0005434180358066555
$ws.Range("P4").Value = 0.0007124942001013113
$ws.Range("Q4").Value = 0.000114833116
$ws.Range("R4").Value = 0.0006889986959999999
$ws.Range("S4").Value = 0.0002406670842940089
$ws.Range("T4").Value = 0.0002467983636304174
# Row 5
$ws.Range("G5").Value = 0.0292345
$ws.Range("H5").Value = 0.058469
$ws.Range("I5").Value = 0.4428765120700495
$ws.Range("J5").Value = 0.346386487911515
$ws.Range("M5").Value = 1.9606995
$ws.Range("N5").Value = 3.921399
$ws.Range("O5").Value = 0.2712524111754306
$ws.Range("P5").Value = 0.2370989514411984
$ws.Range("Q5").Value = 0.05732006953275
$ws.Range("R5").Value = 0.229280278131
$ws.Range("S5").Value = 0.1201313217519656
$ws.Range("T5").Value = 0.08212787307721957
# Row 6
$ws.Range("G6").Value = 0.0292345
$ws.Range("H6").Value = 0.058469
$ws.Range("I6").Value = 0.4428765120700495
$ws.Range("J6").Value = 0.346386487911515
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.02851766666666667
$ws.Range("N6").Value = 0.085553
$ws.Range("O6").Value = 0.003945268433245655
$ws.Range("P6").Value = 0.005172778029639129
$ws.Range("Q6").Value = 0.0008336997261666666
$ws.Range("R6").Value = 0.005002198357
$ws.Range("S6").Value = 0.001747266722895904
$ws.Range("T6").Value = 0.001791780414432544
# Row 7
$ws.Range("G7").Value = 0.0292345
$ws.Range("H7").Value = 0.058469
$ws.Range("I7").Value = 0.4428765120700495
$ws.Range("J7").Value = 0.346386487911515
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.006513000000000001
$ws.Range("N7").Value = 0.019539
$ws.Range("O7").Value = 0.0009010391209798237
$ws.Range("P7").Value = 0.001181383585860448
$ws.Range("Q7").Value = 0.0001904042985
$ws.Range("R7").Value = 0.001142425791
$ws.Range("S7").Value = 0.0003990490631382076
$ws.Range("T7").Value = 0.0004092153111825124
# Row 8
$ws.Range("I8").Value = 0.5571234879299505
$ws.Range("J8").Value = 0.6536135120884849
$ws.Range("M8").Value = 3.1851815
$ws.Range("N8").Value = 6.370363
$ws.Range("O8").Value = 0.4406530230187619
$ws.Range("P8").Value = 0.3851702893788179
$ws.Range("Q8").Value = 0.117138234844
$ws.Range("R8").Value = 0.702829409064
$ws.Range("S8").Value = 0.2454981491510894
$ws.Range("T8").Value = 0.2517525055930273
# Row 9
$ws.Range("I9").Value = 0.5571234879299505
$ws.Range("J9").Value = 0.6536135120884849
$ws.Range("O9").Value = 0.2827048402157753
$ws.Range("P9").Value = 0.3706641033643825
$ws.Range("S9").Value = 0.1575015066356921
$ws.Range("T9").Value = 0.2422710664051232
# Row 10
$ws.Range("I10").Value = 0.5571234879299505
$ws.Range("J10").Value = 0.6536135120884849
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.003928
$ws.Range("N10").Value = 0.011784
$ws.Range("O10").Value = 0.0005434180358066555
$ws.Range("P10").Value = 0.0007124942001013113
$ws.Range("Q10").Value = 0.000144456128
$ws.Range("R10").Value = 0.001300105152
$ws.Range("S10").Value = 0.0003027509515126466
$ws.Range("T10").Value = 0.0004656958364708939
# Row 11
$ws.Range("I11").Value = 0.5571234879299505
$ws.Range("J11").Value = 0.6536135120884849
$ws.Range("M11").Value = 1.9606995
$ws.Range("N11").Value = 3.921399
$ws.Range("O11").Value = 0.2712524111754306
$ws.Range("P11").Value = 0.2370989514411984
$ws.Range("Q11").Value = 0.072106684812
$ws.Range("R11").Value = 0.432640108872
$ws.Range("S11").Value = 0.151121089423465
$ws.Range("T11").Value = 0.1549710783639789
# Row 12
$ws.Range("I12").Value = 0.5571234879299505
$ws.Range("J12").Value = 0.6536135120884849
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.02851766666666667
$ws.Range("N12").Value = 0.085553
$ws.Range("O12").Value = 0.003945268433245655
$ws.Range("P12").Value = 0.005172778029639129
$ws.Range("Q12").Value = 0.001048765709333333
$ws.Range("R12").Value = 0.009438891384
$ws.Range("S12").Value = 0.00219800171034975
$ws.Range("T12").Value = 0.003380997615206584
# Row 13
$ws.Range("I13").Value = 0.5571234879299505
$ws.Range("J13").Value = 0.6536135120884849
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.006513000000000001
$ws.Range("N13").Value = 0.019539
$ws.Range("O13").Value = 0.0009010391209798237
$ws.Range("P13").Value = 0.0007124942001013113
$ws.Range("Q13").Value = 0.000239522088
$ws.Range("R13").Value = 0.002155698792
$ws.Range("S13").Value = 0.0005019900578416159
$ws.Range("T13").Value = 0.0007721682746779358

Write-Host "Applied all cell updates"